# Update the "list_offers" sheet with the current batch of matching offers.
# The product list is fully replaced: some previous rows are dropped (no longer
# within the desired price range), some are kept (possibly with refreshed price
# or link), and several new offers are appended.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet already has 8 data rows (2-9) using two alternating row styles:
#   even rows -> same style as row 2
#   odd rows  -> same style as row 3
# The new data needs 12 rows (2-13), so first propagate the existing alternating
# formatting onto the 4 extra rows (10-13) before the values are written, by
# copying only the formats (not the old contents) of row 2 / row 3 into them.
$ws.Range("A2:C2").Copy() | Out-Null
$ws.Range("A10:C10").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$ws.Range("A12:C12").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$ws.Range("A3:C3").Copy() | Out-Null
$ws.Range("A11:C11").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$ws.Range("A13:C13").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$excel.CutCopyMode = 0

# Write the refreshed "Produto" / "Preço" / "Link" values for every data row.
$ws.Range("A2").Value = 'usado: iphone 12 64gb azul bom - trocafone - apple'
$ws.Range("B2").Value = 3483.92
$ws.Range("C2").Value = 'https://www.magazineluiza.com.br/usado-iphone-12-64gb-azul-bom-trocafone-apple/p/djccka1jka/te/ip12?&seller_id=trocafone'

$ws.Range("A3").Value = 'celular iphone 12 5g 64gb azul - open box'
$ws.Range("B3").Value = 3305.22
$ws.Range("C3").Value = 'https://www.maisbaratofone.com.br/produto/celular-apple-iphone-12-5g/?attribute_pa_condicao=open-box&attribute_pa_armazenamento=64gb&attribute_pa_cor=azul'

$ws.Range("A4").Value = 'iphone 12 64gb branco de vitrine tela 6,1&quot; 4g câmera traseira 12mp+12mp ...'
$ws.Range("B4").Value = 3349
$ws.Range("C4").Value = 'https://www.carrefour.com.br/iphone-12-64gb-branco-de-vitrine-tela-61quot-4g-camera-traseira-12mp12mp-vitrine-mp932521352/p'

$ws.Range("A5").Value = 'celular apple iphone 12 black 64gb vitrine/seminovo com carrregador e cabo'
$ws.Range("B5").Value = 3379.9
$ws.Range("C5").Value = 'https://www.carrefour.com.br/celular-apple-iphone-12-black-64gb-vitrineseminovo-com-carrregador-e-cabo-mp934027523/p'

$ws.Range("A6").Value = 'apple iphone 12 preto 64gb'
$ws.Range("B6").Value = 3130
$ws.Range("C6").Value = 'https://doji.com.br/product/apple-iphone-12-preto-64gb-como-novo'

$ws.Range("A7").Value = 'iphone 12 64gb | celular apple | usado'
$ws.Range("B7").Value = 3015
$ws.Range("C7").Value = 'https://www.enjoei.com.br/p/iphone-12-64gb-87683383?g_campaign=google_shopping'

$ws.Range("A8").Value = 'celular apple iphone 12 white 64gb vitrine/seminovo + acessorios'
$ws.Range("B8").Value = 3499
$ws.Range("C8").Value = 'https://www.carrefour.com.br/celular-apple-iphone-12-white-64gb-vitrineseminovo-acessorios-mp934027474/p'

$ws.Range("A9").Value = 'apple iphone 12 64gb 5g - 12mp ios - tela super retina xdr oled 6.1" - preto'
$ws.Range("B9").Value = 3496
$ws.Range("C9").Value = 'https://www.horizonplay.com.br/apple/iphone/apple-iphone-12-64gb-azul-novo-lacrado-tela-super-retina-xdr-oled-6-1?variant_id=21019&parceiro=8926&srsltid=AfmBOoozsDDR7IV2lzsIi-H9hEmjiMHREVqc1qljiZcBUvg6fMHPQX9JPnE'

$ws.Range("A10").Value = 'aparelho iphone 12 preto 64gb apple seminovo/vitrine sem riscos com acessorios'
$ws.Range("B10").Value = 3289.9
$ws.Range("C10").Value = 'https://www.carrefour.com.br/aparelho-iphone-12-preto-64gb-apple-seminovovitrine-sem-riscos-com-acessorios-mp934027394/p'

$ws.Range("A11").Value = 'iphone 11 apple 64gb e 128gb preto 6,1” 12mp ios (64gb)'
$ws.Range("B11").Value = 3399
$ws.Range("C11").Value = 'https://lumixpel.lojavirtualnuvem.com.br/produtos/iphone-11-apple-64gb-e-128gb-preto-61-12mp-ios/?variant=660570516&pf=mc&srsltid=AfmBOor-FQ3Bwl67FXC1WlzVEsKcidTQI9oQc4rijmidm1cBfqeMVbo2Jns'

$ws.Range("A12").Value = 'smartphone apple iphone 12 64gb câmera dupla'
$ws.Range("B12").Value = 3199
$ws.Range("C12").Value = 'https://www.buscape.com.br/celular/smartphone-apple-iphone-12-64gb-ios?_lc=88&searchterm=iphone%2012%2064gb'

$ws.Range("A13").Value = 'placa de video nvidia geforce rtx 3060 ti 8 gb gddr6 192 bits asus dual-rtx3060ti-o8g-v2'
$ws.Range("B13").Value = 4108.27
$ws.Range("C13").Value = 'https://www.buscape.com.br/placa-de-video/placa-de-video-nvidia-geforce-rtx-3060-ti-8-gb-gddr6-192-bits-asus-dual-rtx3060ti-o8g-v2?_lc=88&searchterm=rtx%203060'
